$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header in H1, copying the header formatting/style
# from the existing "sum" header (G1) so it matches the other header cells.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$ws.Range("H1").Value = "Save"

# Fill in the new "Save" column values for rows 2-45.
$values = @(0,0,0,0,0,0,0,0,0,0,1,1,1,0,1,1,0,1,0,0,1,0,0,1,0,0,0,1,0,0,1,0,0,0,0,0,0,1,1,0,1,0,1,0)

$arr = New-Object 'object[,]' 44,1
for ($i = 0; $i -lt 44; $i++) {
    $arr[$i,0] = $values[$i]
}
$ws.Range("H2:H45").Value = $arr
